# Update "想去人数" (want-to-go count) values in F column across sheets
# 展览 (Exhibitions), 本地生活 (Local life), 全部类型 (All types)

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (rows 5-49) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    5  = 359
    6  = 153
    7  = 153
    8  = 753
    9  = 4117
    14 = 5890
    15 = 459
    16 = 2283
    18 = 156
    19 = 443
    20 = 8853
    22 = 1598
    24 = 2276
    25 = 2363
    27 = 219
    28 = 1921
    29 = 37
    30 = 51
    31 = 318
    33 = 33
    35 = 38
    36 = 16
    37 = 28
    38 = 1214
    39 = 1207
    43 = 1486
    44 = 2387
    46 = 901
    47 = 281
    48 = 1243
    49 = 22
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet: 本地生活 (rows 2-4) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$sheet3Updates = @{
    2 = 675
    3 = 869
    4 = 93
}
foreach ($row in $sheet3Updates.Keys) {
    $ws3.Range("F$row").Value = $sheet3Updates[$row]
}

# --- Sheet: 全部类型 (rows 4-50) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    4  = 675
    5  = 869
    6  = 93
    7  = 359
    9  = 153
    11 = 153
    12 = 753
    13 = 4117
    14 = 4117
    18 = 5890
    19 = 459
    20 = 2283
    22 = 156
    23 = 443
    24 = 8853
    27 = 1598
    28 = 2276
    29 = 2363
    31 = 219
    32 = 1921
    33 = 37
    34 = 51
    35 = 318
    37 = 38
    38 = 28
    39 = 1214
    42 = 1486
    43 = 2387
    44 = 901
    46 = 281
    50 = 1243
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
